$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price values in column D (Price), preserving
# them as text (matching the source data which stores prices as strings).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "271.94"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.04"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.361"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06302"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.660"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.405"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8346"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1627"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08412"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03476"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03145"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09319"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.940"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001708"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04856"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006250"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005475"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001089"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.736"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.321"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01387"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3380"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04688"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006913"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1176"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003600"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01251"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006257"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7888"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1165"
